# Deploying to gh-pages: refresh member roster (process_avatar/final.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Status corrections: several members moved from "在读" (enrolled) to "离开" (left) ---
$ws.Range("F97").Value  = "离开"
$ws.Range("F145").Value = "离开"
$ws.Range("F192").Value = "离开"
$ws.Range("F227").Value = "离开"
$ws.Range("F238").Value = "离开"
$ws.Range("F239").Value = "离开"

# --- Degree corrections in column D (U -> UM/UP, M -> UM) ---
$ws.Range("D188").Value = "UP"
$ws.Range("D189").Value = "UM"
$ws.Range("D190").Value = "UM"
$ws.Range("D214").Value = "UM"
$ws.Range("D218").Value = "UM"
$ws.Range("D219").Value = "UM"
$ws.Range("D222").Value = "UM"

# --- New member row appended at the bottom of the roster ---
$ws.Range("A274").Value = "马文杰"
$ws.Range("B274").Value = "Wenjie Ma"
$ws.Range("D274").Value = "U"
$ws.Range("E274").Value = "/assets/img/members/student/马文杰.jpg"
$ws.Range("F274").Value = "离开"

# --- Keep the autofilter / filter-database range in sync with the new data extent ---
$ws.AutoFilterMode = $false
$ws.Range("A1:F274").AutoFilter()

# The "_xlnm._FilterDatabase" defined name normally tracks the AutoFilter
# range automatically in Excel; make sure it is widened to match too.
$filterName = $wb.Names.Item("Sheet1!_FilterDatabase")
$filterName.RefersTo = '=Sheet1!$A$1:$F$274'

# --- Refresh the saved view state (scroll position / active selection) ---
$excel.ActiveWindow.ScrollRow = 259
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E274").Select()
